{"js": "// Update the date heading and the 25 division problems in the table to the\n// new day's values, per the commit's regenerated output.\n\n// 1) Update the date paragraph above the table.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\ndatePara.insertText(\"2025-10-11 Saturday\", \"Replace\");\n\n// 2) Update every division-problem cell in the table, in row-major order,\n// matching the document's original row/column order.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row -> [oldText -> newText] mapping applied in column order for each\n// non-blank row (rows 0, 4, 8, 12, 16 contain the problems).\nconst rowUpdates = [\n  [0, [\"22\u00f76=\", \"29\u00f76=\", \"48\u00f77=\", \"21\u00f72=\", \"69\u00f72=\"]],\n  [4, [\"55\u00f79=\", \"66\u00f73=\", \"52\u00f73=\", \"21\u00f72=\", \"67\u00f73=\"]],\n  [8, [\"28\u00f74=\", \"87\u00f76=\", \"51\u00f76=\", \"40\u00f73=\", \"68\u00f72=\"]],\n  [12, [\"99\u00f77=\", \"36\u00f78=\", \"83\u00f75=\", \"48\u00f77=\", \"68\u00f73=\"]],\n  [16, [\"91\u00f73=\", \"30\u00f72=\", \"78\u00f77=\", \"99\u00f72=\", \"16\u00f78=\"]],\n];\n\nfor (const [rowIndex, values] of rowUpdates) {\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(rowIndex, col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division problems in the table to the\n# new day's values, per the commit's regenerated output. Every \"old\" string\n# below is unique in the document, so a simple Find/ReplaceAll per pair is\n# unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old = \"2025-10-10 Friday\"; new = \"2025-10-11 Saturday\"},\n    @{old = \"68\u00f76=\";  new = \"22\u00f76=\"},\n    @{old = \"95\u00f75=\";  new = \"29\u00f76=\"},\n    @{old = \"42\u00f77=\";  new = \"48\u00f77=\"},\n    @{old = \"75\u00f72=\";  new = \"21\u00f72=\"},\n    @{old = \"83\u00f78=\";  new = \"69\u00f72=\"},\n    @{old = \"35\u00f74=\";  new = \"55\u00f79=\"},\n    @{old = \"41\u00f77=\";  new = \"66\u00f73=\"},\n    @{old = \"60\u00f75=\";  new = \"52\u00f73=\"},\n    @{old = \"31\u00f78=\";  new = \"21\u00f72=\"},\n    @{old = \"21\u00f74=\";  new = \"67\u00f73=\"},\n    @{old = \"16\u00f73=\";  new = \"28\u00f74=\"},\n    @{old = \"29\u00f75=\";  new = \"87\u00f76=\"},\n    @{old = \"98\u00f77=\";  new = \"51\u00f76=\"},\n    @{old = \"72\u00f79=\";  new = \"40\u00f73=\"},\n    @{old = \"84\u00f77=\";  new = \"68\u00f72=\"},\n    @{old = \"62\u00f78=\";  new = \"99\u00f77=\"},\n    @{old = \"74\u00f75=\";  new = \"36\u00f78=\"},\n    @{old = \"60\u00f77=\";  new = \"83\u00f75=\"},\n    @{old = \"26\u00f74=\";  new = \"48\u00f77=\"},\n    @{old = \"93\u00f77=\";  new = \"68\u00f73=\"},\n    @{old = \"67\u00f72=\";  new = \"91\u00f73=\"},\n    @{old = \"90\u00f76=\";  new = \"30\u00f72=\"},\n    @{old = \"51\u00f79=\";  new = \"78\u00f77=\"},\n    @{old = \"78\u00f79=\";  new = \"99\u00f72=\"},\n    @{old = \"15\u00f75=\";  new = \"16\u00f78=\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute($r.old, $false, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)\n}\n"}
